$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "XML"
$ws.Range("B4").Value = "Book"
$ws.Range("C4").Value = "2.0.0"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 4

$ws.Range("A5").Value = "JSON"
$ws.Range("B5").Value = "GetPassengerList"
$ws.Range("C5").Value = "3.4.1"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 5

$ws.Range("A6").Value = "JSON"
$ws.Range("B6").Value = "GetFlightList"
$ws.Range("C6").Value = "2.3.4"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 3

$ws.Range("A7").Select()
